$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.705.10'
$ws.Range("E2").Value = '  -0.54%  '
$ws.Range("D3").Value = '1.628.96'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.39'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("E9").Value = '  -0.99%  '
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.26'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.854.83'
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.625.73'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.554'
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '25.722.42'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.43'
$ws.Range("E21").Value = '  -1.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.91'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("E24").Value = '  -0.84%  '
$ws.Range("E25").Value = '  +1.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.94'
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("E30").Value = '  -0.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0492'
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("E34").Value = '  -0.72%  '
$ws.Range("E35").Value = '  -0.74%  '
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").Value = '1.142.34'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.541'
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("E42").Value = '  -0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.70'
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.53'
$ws.Range("E44").Value = '  -0.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.805'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("D46").Value = '1.764.41'
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.06'
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("E48").Value = '  +0.83%  '
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.37'
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.43'
$ws.Range("E51").Value = '  +4.58%  '
